# edit.ps1
# Applies the diff: fixes row 170 (A170/D170 text->number), appends new
# weather-data rows 171-223 (normal typed rows), and appends row 224 whose
# Temperature (A) and Current Hour (D) values stay as TEXT (matching the
# source data's "not yet normalized" last row, same as old row 170 before
# the fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 170: convert A170 and D170 from text to real numbers.
# B170 ("Mostly Cloudy") and C170 ("08/14/2024") are left untouched.
$ws.Cells.Item(170, 1).Value = 31
$ws.Cells.Item(170, 4).Value = 17

# --- New rows 171-223: Temperature(A, number), Weather Condition(B, text),
# Date Today(C, text), Current Hour(D, number)
$rows171to223 = @(
    @(171, 29, "Cloudy", "08/15/2024", 8),
    @(172, 29, "Cloudy", "08/15/2024", 8),
    @(173, 29, "Cloudy", "08/15/2024", 8),
    @(174, 29, "Cloudy", "08/15/2024", 8),
    @(175, 29, "Cloudy", "08/15/2024", 8),
    @(176, 30, "Cloudy", "08/15/2024", 9),
    @(177, 30, "Cloudy", "08/15/2024", 9),
    @(178, 30, "Cloudy", "08/15/2024", 9),
    @(179, 30, "Cloudy", "08/15/2024", 9),
    @(180, 31, "Cloudy", "08/15/2024", 9),
    @(181, 31, "Cloudy", "08/15/2024", 9),
    @(182, 31, "Cloudy", "08/15/2024", 10),
    @(183, 31, "Cloudy", "08/15/2024", 10),
    @(184, 31, "Cloudy", "08/15/2024", 10),
    @(185, 31, "Cloudy", "08/15/2024", 10),
    @(186, 31, "Cloudy", "08/15/2024", 10),
    @(187, 32, "Cloudy", "08/15/2024", 10),
    @(188, 32, "Cloudy", "08/15/2024", 11),
    @(189, 32, "Cloudy", "08/15/2024", 11),
    @(190, 32, "Cloudy", "08/15/2024", 11),
    @(191, 32, "Cloudy", "08/15/2024", 11),
    @(192, 32, "Cloudy", "08/15/2024", 11),
    @(193, 32, "Cloudy", "08/15/2024", 12),
    @(194, 32, "Cloudy", "08/15/2024", 12),
    @(195, 32, "Cloudy", "08/15/2024", 12),
    @(196, 32, "Cloudy", "08/15/2024", 12),
    @(197, 32, "Cloudy", "08/15/2024", 12),
    @(198, 32, "Cloudy", "08/15/2024", 12),
    @(199, 32, "Cloudy", "08/15/2024", 13),
    @(200, 32, "Mostly Cloudy", "08/15/2024", 13),
    @(201, 32, "Mostly Cloudy", "08/15/2024", 13),
    @(202, 32, "Mostly Cloudy", "08/15/2024", 13),
    @(203, 32, "Mostly Cloudy", "08/15/2024", 13),
    @(204, 32, "Mostly Cloudy", "08/15/2024", 13),
    @(205, 32, "Mostly Cloudy", "08/15/2024", 14),
    @(206, 32, "Cloudy", "08/15/2024", 14),
    @(207, 32, "Cloudy", "08/15/2024", 14),
    @(208, 32, "Cloudy", "08/15/2024", 14),
    @(209, 32, "Cloudy", "08/15/2024", 14),
    @(210, 32, "Cloudy", "08/15/2024", 14),
    @(211, 32, "Cloudy", "08/15/2024", 15),
    @(212, 32, "Mostly Cloudy", "08/15/2024", 15),
    @(213, 32, "Mostly Cloudy", "08/15/2024", 15),
    @(214, 32, "Cloudy", "08/15/2024", 15),
    @(215, 32, "Cloudy", "08/15/2024", 15),
    @(216, 32, "Cloudy", "08/15/2024", 15),
    @(217, 31, "Cloudy", "08/15/2024", 16),
    @(218, 31, "Cloudy", "08/15/2024", 16),
    @(219, 31, "Cloudy", "08/15/2024", 16),
    @(220, 31, "Cloudy", "08/15/2024", 16),
    @(221, 31, "Cloudy", "08/15/2024", 16),
    @(222, 31, "Cloudy", "08/15/2024", 17),
    @(223, 31, "Cloudy", "08/15/2024", 17)
)

foreach ($row in $rows171to223) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $row[3]
    $cCell.Style = "Normal"

    $ws.Cells.Item($r, 4).Value = $row[4]
}

# --- New row 224: all four values stored as TEXT (A and D are numeric-looking
# strings "30"/"17" that must NOT be auto-converted to numbers or dates).
$aCell = $ws.Cells.Item(224, 1)
$aCell.NumberFormat = "@"
$aCell.Value = "30"
$aCell.Style = "Normal"

$ws.Cells.Item(224, 2).Value = "Cloudy"

$cCell224 = $ws.Cells.Item(224, 3)
$cCell224.NumberFormat = "@"
$cCell224.Value = "08/15/2024"
$cCell224.Style = "Normal"

$dCell = $ws.Cells.Item(224, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "17"
$dCell.Style = "Normal"
